# #5: property boat&car done
#
# The "汽車" (cars) sheet originally only had an (undeclared/ad-hoc) header
# row whose cells actually held copies of row-2's data instead of real
# column labels, and it only carried columns A:G (index, name, capacity,
# owner, register_date, register_reason, acquire_value).
#
# This change:
#   1. Rewrites row 1 into a proper header row (name, capacity, owner,
#      register_date, register_reason, acquire_value, property_category,
#      category, date, legislator_name, legislator_id, source_file, index)
#      matching the header convention used by every other sheet in the
#      workbook, carrying over the existing bold/bordered/centered header
#      style.
#   2. Extends each data row (2-5) with the property_category / category /
#      date / legislator_name / legislator_id / source_file / index
#      columns (H:N) that every other property-type sheet already has.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- 1. Rebuild the header row -------------------------------------------
# Column A's header-style cells (e.g. A2) already carry the bold + thin
# border + centered alignment used for every header in this workbook, so
# reuse that exact formatting for the whole new header row B1:N1 rather
# than re-deriving it property by property.
$ws.Range("A2").Copy()
$ws.Range("B1:N1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- 2. Extend data rows 2-5 with columns H:N -----------------------------
$carRows = @(52, 53, 54, 55)

for ($i = 0; $i -lt $carRows.Length; $i++) {
    $r = $i + 2
    $ws.Range("H$r").Value = "land"
    $ws.Range("I$r").Value = "normal"
    # Leading apostrophe forces text so "2012-04-16" isn't re-interpreted
    # as a date serial (matches every other sheet, where this column is
    # stored as the literal text "2012-04-16").
    $ws.Range("J$r").Value = "'2012-04-16"
    $ws.Range("K$r").Value = "紀國棟"
    $ws.Range("L$r").Value = 918
    $ws.Range("M$r").Value = "tmpf6b91"
    $ws.Range("N$r").Value = $carRows[$i]
}
